$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Prix Spot")

# Insert a new column before DP, shifting DP:ET to the right (to DQ:EU)
$ws.Range("DP1:DP25").EntireColumn.Insert()

# Set header value for new column DP1 (new date column "11-nov")
$ws.Range("DP1").Value = "11-nov"

# Set data cells DP2:DP25 to "-" (no data for this date yet)
$ws.Range("DP2:DP25").Value = "-"
